# "Test run complete, bugfixes needed" -- refresh the simulated mean flow
# rate results (column B, rows 2-9) with the new test-run values, and
# resize the data columns to fit the refreshed content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 54.473378555443993
$ws.Range("B3").Value = 32.598269851519511
$ws.Range("B4").Value = 13.411851961272751
$ws.Range("B5").Value = 55.870498547441031
$ws.Range("B6").Value = 33.73868240861232
$ws.Range("B7").Value = 14.639264531641869
$ws.Range("B8").Value = 38.083375505035008
$ws.Range("B9").Value = 15.681388231401629

$ws.Columns.Item(1).ColumnWidth = 25.666666666666668
$ws.Columns.Item(2).ColumnWidth = 13.666666666666666
$ws.Columns.Item(3).ColumnWidth = 4
